$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.702.58"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.35"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.81"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4676"
$ws.Range("E7").Value = "  +4.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3934"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.91"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08026"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("E11").Value = "  +0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.72"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.48"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.122"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001046"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.53"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06617"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.715.62"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.483"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.98"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.082.99"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.28"
$ws.Range("E27").Value = "  +5.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.12"
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.083"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.51"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9667"
$ws.Range("E32").Value = "  +3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09473"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.444"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.597"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.309"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02259"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06060"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.234"
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.126"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5974"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1890"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.21"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.270"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5693"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.24"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.386"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.932"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06840"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.90"
$ws.Range("E51").Value = "  +5.19%  "
